# Auto-generated edit script: update Leve price/profit columns (H-N) per scheduled price-runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1111291.5
$ws.Range("I2").Value = 1250194.9
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 1250194.9
$ws.Range("L2").Value = 65
$ws.Range("M2").Value = -1250081.9
$ws.Range("N2").Value = -291
$ws.Range("H9").Value = 384839.16
$ws.Range("I9").Value = 3333363.2
$ws.Range("K9").Value = 3333363.2
$ws.Range("M9").Value = -3333194.2
$ws.Range("H18").Value = 3599
$ws.Range("I18").Value = 2498.75
$ws.Range("K18").Value = 2498.75
$ws.Range("M18").Value = -2214.75
$ws.Range("H51").Value = 4199
$ws.Range("J51").Value = 4199
$ws.Range("L51").Value = 4199
$ws.Range("N51").Value = -5167
$ws.Range("H69").Value = 1000000000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 1000000000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 3000000000
$ws.Range("M69").ClearContents()  # was -23123
$ws.Range("N69").Value = -3000001748
$ws.Range("H72").Value = 1000000000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 1000000000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 9000000000
$ws.Range("M72").ClearContents()  # was -67623
$ws.Range("N72").Value = -9000008736
$ws.Range("H98").Value = 7991.433
$ws.Range("I98").Value = 9837.130999999999
$ws.Range("K98").Value = 9837.130999999999
$ws.Range("M98").Value = -8339.130999999999
$ws.Range("H122").Value = 7991.433
$ws.Range("I122").Value = 9837.130999999999
$ws.Range("K122").Value = 29511.393
$ws.Range("M122").Value = -27061.393
$ws.Range("H137").Value = 1890.2
$ws.Range("I137").Value = 1916.2222
$ws.Range("K137").Value = 5748.6666
$ws.Range("M137").Value = -3198.6666
$ws.Range("H141").Value = 10613.214
$ws.Range("I141").Value = 10509.889
$ws.Range("K141").Value = 31529.667
$ws.Range("M141").Value = -26349.667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3186.7827
$ws.Range("I2").Value = 2089.4285
$ws.Range("J2").Value = 4893.778
$ws.Range("K2").Value = 2089.4285
$ws.Range("L2").Value = 4893.778
$ws.Range("M2").Value = -1976.4285
$ws.Range("N2").Value = -5119.778
$ws.Range("H4").Value = 475
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()  # was -10241.5
$ws.Range("H74").Value = 2384.9048
$ws.Range("I74").Value = 2391.0625
$ws.Range("K74").Value = 2391.0625
$ws.Range("M74").Value = -1517.0625
$ws.Range("H77").Value = 2384.9048
$ws.Range("I77").Value = 2391.0625
$ws.Range("K77").Value = 11955.3125
$ws.Range("M77").Value = -7587.3125
$ws.Range("H95").Value = 64851.832
$ws.Range("J95").Value = 64851.832
$ws.Range("L95").Value = 64851.832
$ws.Range("N95").Value = -70343.83199999999
$ws.Range("H116").Value = 3186.7827
$ws.Range("I116").Value = 2089.4285
$ws.Range("J116").Value = 4893.778
$ws.Range("K116").Value = 2089.4285
$ws.Range("L116").Value = 4893.778
$ws.Range("M116").Value = 204.5715
$ws.Range("N116").Value = -9481.778

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3186.7827
$ws.Range("I3").Value = 2089.4285
$ws.Range("J3").Value = 4893.778
$ws.Range("K3").Value = 2089.4285
$ws.Range("L3").Value = 4893.778
$ws.Range("M3").Value = -1975.4285
$ws.Range("N3").Value = -5121.778
$ws.Range("H22").Value = 10328.8
$ws.Range("I22").Value = 12842.25
$ws.Range("J22").Value = 275
$ws.Range("K22").Value = 12842.25
$ws.Range("L22").Value = 275
$ws.Range("M22").Value = -12669.25
$ws.Range("N22").Value = -621
$ws.Range("H57").Value = 74500
$ws.Range("J57").Value = 74500
$ws.Range("L57").Value = 74500
$ws.Range("N57").Value = -75940
$ws.Range("H60").Value = 35000
$ws.Range("J60").Value = 35000
$ws.Range("L60").Value = 35000
$ws.Range("N60").Value = -36198
$ws.Range("H136").Value = 74500
$ws.Range("J136").Value = 74500
$ws.Range("L136").Value = 74500
$ws.Range("N136").Value = -84700

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2665.6667
$ws.Range("J22").Value = 2665.6667
$ws.Range("L22").Value = 2665.6667
$ws.Range("N22").Value = -3365.6667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1300.6
$ws.Range("I34").Value = 501
$ws.Range("J34").Value = 1500.5
$ws.Range("K34").Value = 1503
$ws.Range("L34").Value = 4501.5
$ws.Range("M34").Value = -1419
$ws.Range("N34").Value = -4669.5
$ws.Range("H39").Value = 1485.2858
$ws.Range("I39").Value = 1099.25
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 3297.75
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = -3003.75
$ws.Range("N39").Value = -6588
$ws.Range("H55").Value = 357440.44
$ws.Range("I55").Value = 454651.47
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 1363954.41
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -1363777.41
$ws.Range("N55").Value = -3354
$ws.Range("H131").Value = 3337595.5
$ws.Range("J131").Value = 3708272.8
$ws.Range("L131").Value = 11124818.4
$ws.Range("N131").Value = -11134898.4
$ws.Range("H137").Value = 14030.5
$ws.Range("J137").Value = 14128.0625
$ws.Range("L137").Value = 42384.1875
$ws.Range("N137").Value = -52584.1875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4466.196
$ws.Range("I80").Value = 3022.4194
$ws.Range("K80").Value = 3022.4194
$ws.Range("M80").Value = -2024.4194
$ws.Range("H83").Value = 4466.196
$ws.Range("I83").Value = 3022.4194
$ws.Range("K83").Value = 15112.097
$ws.Range("M83").Value = -10120.097
$ws.Range("H139").Value = 102983
$ws.Range("J139").Value = 102983
$ws.Range("L139").Value = 102983
$ws.Range("N139").Value = -113263

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3532.4443
$ws.Range("I22").Value = 2399
$ws.Range("J22").Value = 3674.125
$ws.Range("K22").Value = 2399
$ws.Range("L22").Value = 3674.125
$ws.Range("M22").Value = -2104
$ws.Range("N22").Value = -4264.125
$ws.Range("H27").Value = 3532.4443
$ws.Range("I27").Value = 2399
$ws.Range("J27").Value = 3674.125
$ws.Range("K27").Value = 2399
$ws.Range("L27").Value = 3674.125
$ws.Range("M27").Value = -2292
$ws.Range("N27").Value = -3888.125
$ws.Range("H76").Value = 40896
$ws.Range("J76").Value = 40896
$ws.Range("L76").Value = 40896
$ws.Range("N76").Value = -41572
$ws.Range("H79").Value = 40896
$ws.Range("J79").Value = 40896
$ws.Range("L79").Value = 40896
$ws.Range("N79").Value = -43236
$ws.Range("H132").Value = 3523.3076
$ws.Range("I132").Value = 3143.6875
$ws.Range("J132").Value = 4130.7
$ws.Range("K132").Value = 9431.0625
$ws.Range("L132").Value = 12392.1
$ws.Range("M132").Value = -6901.0625
$ws.Range("N132").Value = -17452.1
$ws.Range("H137").Value = 87616.39999999999
$ws.Range("J137").Value = 84360.664
$ws.Range("L137").Value = 84360.664
$ws.Range("N137").Value = -94560.664
$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2308.7646
$ws.Range("I81").Value = 1537.3334
$ws.Range("K81").Value = 3074.6668
$ws.Range("M81").Value = -2013.6668
$ws.Range("H84").Value = 2308.7646
$ws.Range("I84").Value = 1537.3334
$ws.Range("K84").Value = 15373.334
$ws.Range("M84").Value = -10069.334
$ws.Range("H126").Value = 3205.5557
$ws.Range("I126").Value = 2505.7144
$ws.Range("J126").Value = 5655
$ws.Range("K126").Value = 7517.1432
$ws.Range("L126").Value = 16965
$ws.Range("M126").Value = -5047.1432
$ws.Range("N126").Value = -21905
$ws.Range("H128").Value = 89800
$ws.Range("J128").Value = 89800
$ws.Range("L128").Value = 89800
$ws.Range("N128").Value = -99760

